$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.409.16"
$ws.Range("E2").Value = "  -3.34%  "
$ws.Range("D3").Value = "3.689.30"
$ws.Range("E3").Value = "  -3.58%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.13%  "
$ws.Range("D7").Value = "3.684.98"
$ws.Range("E7").Value = "  -3.72%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.24%  "
$ws.Range("E12").Value = "  -3.39%  "
$ws.Range("E13").Value = "  -6.41%  "
$ws.Range("E14").Value = "  -5.51%  "
$ws.Range("D15").Value = "4.302.11"
$ws.Range("E15").Value = "  -3.52%  "
$ws.Range("D16").Value = "3.687.71"
$ws.Range("E16").Value = "  -3.51%  "
$ws.Range("D17").Value = "67.459.54"
$ws.Range("E17").Value = "  -3.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.18%  "
$ws.Range("E19").Value = "  -4.27%  "
$ws.Range("E20").Value = "  -3.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "490.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.723"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("E25").Value = "  -6.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000138"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.97%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.69%  "
$ws.Range("E31").Value = "  -7.07%  "
$ws.Range("E32").Value = "  -4.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("D34").Value = "3.826.46"
$ws.Range("E34").Value = "  -3.49%  "
$ws.Range("E35").Value = "  -5.09%  "
$ws.Range("D36").Value = "3.627.07"
$ws.Range("E36").Value = "  -3.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.992"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.131"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.322"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "434.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.30%  "
$ws.Range("E44").Value = "  -5.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.85%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.28%  "
$ws.Range("D50").Value = "2.752.85"
$ws.Range("E50").Value = "  -5.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0347"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.86%  "
